# CANmailboxUsage.xlsx — "Fix input signals and get temp diodes working"
#
# 1. "CAN IDs" sheet: add the 141 / "Board Temp" readings in row 9 (F9/H9),
#    and remove the now-unused BP1/BP2/IO1-6 lookup table that lived in
#    columns P:Q (rows 16-24).
# 2. "Mailboxes (CANB)" sheet: add a 15th mailbox row (row 17) with the same
#    Output / CPU1 / HO_CAN->CAN2->analogInputs pattern used by the rows
#    above it (rows 12-16, mailboxes 9-13).
# 3. Restore the on-screen selections left behind by the author's edit.

$wb = $excel.ActiveWorkbook

# --- "CAN IDs" sheet -------------------------------------------------
$wsCanIds = $wb.Worksheets.Item("CAN IDs")
$wsCanIds.Activate() | Out-Null

$wsCanIds.Range("F9").Value = 141
$wsCanIds.Range("H9").Value = "Board Temp"

# Drop the old BP1/BP2/IO1-6 helper table (columns P:Q, rows 16-24)
$wsCanIds.Range("P16:Q24").Clear() | Out-Null

$wsCanIds.Range("J23").Select() | Out-Null

# --- "Mailboxes (CANB)" sheet -----------------------------------------
$wsCanB = $wb.Worksheets.Item("Mailboxes (CANB)")
$wsCanB.Activate() | Out-Null

$wsCanB.Range("C17").Value = "Output"
$wsCanB.Range("D17").Value = "CPU1"
$wsCanB.Range("E17").Value = "HO_CAN->CAN2->analogInputs"

$wsCanB.Range("C17:E17").Select() | Out-Null

# Try to restore the saved window position (best effort; some hosts don't
# persist this purely cosmetic piece of UI state).
$excel.ActiveWindow.Left = -120
